$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("historical_data")

$years = 2000..2023

# ---- Section 1: EMBER Capacity (GW) ----
$ws.Cells.Item(28, 1).Value = "EMBER Capacity (GW)"

$ws.Cells.Item(29, 1).Value = "model_fuel"
for ($i = 0; $i -lt $years.Count; $i++) {
    $ws.Cells.Item(29, 2 + $i).Value = $years[$i]
}

$capacityData = @{
    "bioenergy" = @(0.05,0.05,0.05,0.01,0.01,0.01,0.01,0.01,0.01,0.01,0.01,0.01,0.01,0.03,0.04,0.05,0.06,0.05,0.07,0.06,0.05,0.05,0.05,0.05)
    "coal"      = @(5.62,5.62,5.62,5.62,5.5,5.5,5.5,5.5,5.5,5.63,5.63,6.32,6.32,6.11,5.9,5.0599999999999996,5.0599999999999996,5.1100000000000003,5.1100000000000003,5.1100000000000003,5.1100000000000003,5.1100000000000003,5.1100000000000003,5.1100000000000003)
    "gas"       = @(1.1100000000000001,1.1100000000000001,1.1100000000000001,1.1100000000000001,1.1100000000000001,1.1100000000000001,1.1100000000000001,1.1499999999999999,1.1499999999999999,1.1499999999999999,1.1499999999999999,1.2,1.2,1.2,1.2,1.2,1.2,1.2,1.2,1.2,1.2,1.2,1.2,1.2)
    "hydro"     = @(1.02,0.84,1.08,1.65,1.98,1.98,1.98,2.0099999999999998,2.12,2.14,2.1800000000000002,2.2400000000000002,2.3199999999999998,2.34,2.36,2.36,2.36,2.5099999999999998,2.5099999999999998,2.5099999999999998,2.5099999999999998,2.5099999999999998,2.5299999999999998,2.5299999999999998)
    "nuclear"   = @(3.53,3.53,2.72,2.72,2.72,2.72,2.72,1.89,1.89,1.89,1.89,1.89,1.91,1.98,1.98,1.98,1.97,1.97,2.0099999999999998,2.0099999999999998,2.0099999999999998,2.0099999999999998,2.0099999999999998,2.0099999999999998)
    "oil"       = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    "solar"     = @(0,0,0,0,0,0,0,0,0,0,0.03,0.15,0.92,1.04,1.03,1.03,1.03,1.03,1.03,1.04,1.1000000000000001,1.27,1.74,2.94)
    "wind"      = @(0,0,0,0,0,0.01,0.03,0.03,0.11,0.33,0.49,0.54,0.68,0.68,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7,0.7)
}

$capacityRows = @("bioenergy","coal","gas","hydro","nuclear","oil","solar","wind")
$r = 30
foreach ($fuel in $capacityRows) {
    $ws.Cells.Item($r, 1).Value = $fuel
    $vals = $capacityData[$fuel]
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
    $r++
}

# ---- Section 2: EMBER Generation (TWh) ----
$ws.Cells.Item(40, 1).Value = "EMBER Generation (TWh)"

$ws.Cells.Item(41, 1).Value = "model_fuel"
for ($i = 0; $i -lt $years.Count; $i++) {
    $ws.Cells.Item(41, 2 + $i).Value = $years[$i]
}

$generationData = @{
    "bioenergy" = @(0,0,0,0,0,0,0,0,0.02,0.01,0.04,0.06,0.07,0.11,0.2,0.27,0.35,0.4,1.57,1.82,1.7,2.59,2.2400000000000002,2.1800000000000002)
    "coal"      = @(16.940000000000001,19.5,17.18,19.239999999999998,18.899999999999999,18.46,19.09,22.37,23.18,21.1,22.61,27.53,22.87,19.39,21.31,22.53,19.37,20.92,18.66,17.2,13.51,17.09,21.79,11.56)
    "gas"       = @(1.91,1.91,1.54,1.76,1.49,1.73,2.16,2.34,2.36,1.96,1.97,2.08,2.36,2.34,2.14,1.86,2.0499999999999998,1.93,2.02,2.15,2.29,3.05,2.0499999999999998,1.56)
    "hydro"     = @(2.63,1.65,2.12,2.99,3.14,4.3,4.21,2.83,2.79,3.43,5.03,2.87,3.18,4.04,4.5999999999999996,5.65,3.88,2.83,5.15,2.93,2.82,4.82,3.8,3.11)
    "nuclear"   = @(18.18,19.55,20.22,17.28,16.82,18.649999999999999,19.489999999999998,14.64,15.77,15.26,15.25,16.309999999999999,15.78,14.17,15.87,15.38,15.78,15.55,16.13,16.559999999999999,16.63,16.489999999999998,16.46,16.16)
    "oil"       = @(0.93,0.83,1.02,1.01,1.01,0.75,0.47,0.59,0.25,0.31,0.34,0.15,0.2,0.21,0.21,0.19,0.28999999999999998,0.31,0.32,0.32,0.24,0.28000000000000003,0.35,0.35)
    "solar"     = @(0,0,0,0,0,0,0,0,0,0,0.01,0.1,0.78,1.39,1.26,1.38,1.39,1.4,1.34,1.42,1.47,1.47,2.09,3.52)
    "wind"      = @(0,0,0,0,0,0,0.02,0.05,0.12,0.24,0.68,0.86,1.22,1.37,1.33,1.45,1.42,1.5,1.32,1.32,1.48,1.43,1.44,1.55)
}

$generationRows = @("bioenergy","coal","gas","hydro","nuclear","oil","solar","wind")
$r = 42
foreach ($fuel in $generationRows) {
    $ws.Cells.Item($r, 1).Value = $fuel
    $vals = $generationData[$fuel]
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
    $r++
}
